$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2").Value = "3/11/2001"
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F8").Select()
